$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7201011180877686
$ws.Range("B1").Value = 3.112647533416748
$ws.Range("C1").Value = 2.885094881057739
$ws.Range("D1").Value = 2.398935317993164
$ws.Range("E1").Value = 2.123542785644531
